# Updated flu processing to utilize error code datatable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new "multiple matches" error code (note: value typed before code, matching
# the original authoring order so the shared-string table lines up)
$ws.Range("B3").Value = "Multiple matches found in MRRS Report. Sailor's info will have to be checked manually."
$ws.Range("A3").Value = "MULTIPLE_RECORDS_FOUND"

# New "Process" column (C) with header + existing rows' process classification
$ws.Range("C1").Value = "Process"
$ws.Range("C2").Value = "Record Search"
$ws.Range("C3").Value = "Record Search"

# Row 4 - influenza overdue
$ws.Range("A4").Value = "INFLUENZA_OVERDUE"
$ws.Range("B4").Value = "Influenza out of date in MRRS. Verify with physical record and MHS Genesis."
$ws.Range("C4").Value = "Influeza Verification"

# Row 5 - no documented influenza
$ws.Range("A5").Value = "NO_DOCUMENTED_INFLUENZA"
$ws.Range("B5").Value = "No documented influenza in MRRS. Verify with physical record and MHS Genesis."
$ws.Range("C5").Value = "Influeza Verification"

# Resize columns to fit the new, longer content (mirrors Excel's "best fit" resize)
$ws.Columns.Item(1).ColumnWidth = 26.5
$ws.Columns.Item(2).ColumnWidth = 73.5
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668

# Final selection left on A4, as in the saved workbook
$ws.Range("A4").Select() | Out-Null
